$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(18, 8).Value = 2100  # ALC!H18: 2116.8333 -> 2100
$ws.Cells.Item(18, 9).Value = 2120  # ALC!I18: 2116.8333 -> 2120
$ws.Cells.Item(18, 10).Value = 2000  # ALC!J18: 0 -> 2000
$ws.Cells.Item(18, 11).Value = 2120  # ALC!K18: 2116.8333 -> 2120
$ws.Cells.Item(18, 12).Value = 2000  # ALC!L18: 0 -> 2000
$ws.Cells.Item(18, 13).Value = -1836  # ALC!M18: -1832.8333 -> -1836
$ws.Cells.Item(18, 14).Value = -2568  # ALC!N18: None -> -2568

$ws.Cells.Item(62, 8).Value = 6322.154  # ALC!H62: 6653.8184 -> 6322.154
$ws.Cells.Item(62, 9).Value = 4480.3335  # ALC!I62: 4721.5 -> 4480.3335
$ws.Cells.Item(62, 10).Value = 6874.7  # ALC!J62: 7083.222 -> 6874.7
$ws.Cells.Item(62, 11).Value = 4480.3335  # ALC!K62: 4721.5 -> 4480.3335
$ws.Cells.Item(62, 12).Value = 6874.7  # ALC!L62: 7083.222 -> 6874.7
$ws.Cells.Item(62, 13).Value = -3856.3335  # ALC!M62: -4097.5 -> -3856.3335
$ws.Cells.Item(62, 14).Value = -8122.7  # ALC!N62: -8331.222 -> -8122.7

$ws.Cells.Item(65, 8).Value = 6322.154  # ALC!H65: 6653.8184 -> 6322.154
$ws.Cells.Item(65, 9).Value = 4480.3335  # ALC!I65: 4721.5 -> 4480.3335
$ws.Cells.Item(65, 10).Value = 6874.7  # ALC!J65: 7083.222 -> 6874.7
$ws.Cells.Item(65, 11).Value = 22401.6675  # ALC!K65: 23607.5 -> 22401.6675
$ws.Cells.Item(65, 12).Value = 34373.5  # ALC!L65: 35416.11 -> 34373.5
$ws.Cells.Item(65, 13).Value = -19281.6675  # ALC!M65: -20487.5 -> -19281.6675
$ws.Cells.Item(65, 14).Value = -40613.5  # ALC!N65: -41656.11 -> -40613.5

$ws.Cells.Item(107, 8).Value = 376.44446  # ALC!H107: 498.85715 -> 376.44446
$ws.Cells.Item(107, 9).Value = 392.3125  # ALC!I107: 518 -> 392.3125
$ws.Cells.Item(107, 10).Value = 249.5  # ALC!J107: 250 -> 249.5
$ws.Cells.Item(107, 11).Value = 392.3125  # ALC!K107: 518 -> 392.3125
$ws.Cells.Item(107, 12).Value = 249.5  # ALC!L107: 250 -> 249.5
$ws.Cells.Item(107, 13).Value = 1527.6875  # ALC!M107: 1402 -> 1527.6875
$ws.Cells.Item(107, 14).Value = -4089.5  # ALC!N107: -4090 -> -4089.5

$ws.Cells.Item(137, 8).Value = 8052.8096  # ALC!H137: 7732.227 -> 8052.8096
$ws.Cells.Item(137, 9).Value = 1675.3  # ALC!I137: 1613.909 -> 1675.3
$ws.Cells.Item(137, 11).Value = 5025.9  # ALC!K137: 4841.727000000001 -> 5025.9
$ws.Cells.Item(137, 13).Value = -2475.9  # ALC!M137: -2291.727000000001 -> -2475.9

$ws.Cells.Item(138, 8).Value = 5043.887  # ALC!H138: 5116.586 -> 5043.887
$ws.Cells.Item(138, 9).Value = 5912.625  # ALC!I138: 6157.9565 -> 5912.625
$ws.Cells.Item(138, 10).Value = 4600.2764  # ALC!J138: 4606.9785 -> 4600.2764
$ws.Cells.Item(138, 11).Value = 17737.875  # ALC!K138: 18473.8695 -> 17737.875
$ws.Cells.Item(138, 12).Value = 13800.8292  # ALC!L138: 13820.9355 -> 13800.8292
$ws.Cells.Item(138, 13).Value = -12597.875  # ALC!M138: -13333.8695 -> -12597.875
$ws.Cells.Item(138, 14).Value = -24080.8292  # ALC!N138: -24100.9355 -> -24080.8292

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 17615.045  # ARM!H32: 18129.354 -> 17615.045
$ws.Cells.Item(32, 9).Value = 7633.919  # ARM!I32: 8018.7144 -> 7633.919
$ws.Cells.Item(32, 11).Value = 7633.919  # ARM!K32: 8018.7144 -> 7633.919
$ws.Cells.Item(32, 13).Value = -7346.919  # ARM!M32: -7731.7144 -> -7346.919

$ws.Cells.Item(61, 8).Value = 2400  # ARM!H61: 2398.3333 -> 2400
$ws.Cells.Item(61, 9).Value = 0  # ARM!I61: 2397.5 -> 0
$ws.Cells.Item(61, 11).Value = 0  # ARM!K61: 2397.5 -> 0
$ws.Cells.Item(61, 13).ClearContents()  # ARM!M61: -2185.5 -> (removed)

$ws.Cells.Item(97, 8).Value = 511.10526  # ARM!H97: 587.5263 -> 511.10526
$ws.Cells.Item(97, 9).Value = 505.85715  # ARM!I97: 575.5333000000001 -> 505.85715
$ws.Cells.Item(97, 10).Value = 525.8  # ARM!J97: 632.5 -> 525.8
$ws.Cells.Item(97, 11).Value = 505.85715  # ARM!K97: 575.5333000000001 -> 505.85715
$ws.Cells.Item(97, 12).Value = 525.8  # ARM!L97: 632.5 -> 525.8
$ws.Cells.Item(97, 13).Value = -9.85714999999999  # ARM!M97: -79.53330000000005 -> -9.85714999999999
$ws.Cells.Item(97, 14).Value = -1517.8  # ARM!N97: -1624.5 -> -1517.8

$ws.Cells.Item(102, 8).Value = 1264.6666  # ARM!H102: 1290.7391 -> 1264.6666
$ws.Cells.Item(102, 9).Value = 1252.6666  # ARM!I102: 1294.6428 -> 1252.6666
$ws.Cells.Item(102, 11).Value = 1252.6666  # ARM!K102: 1294.6428 -> 1252.6666
$ws.Cells.Item(102, 13).Value = 369.3334  # ARM!M102: 327.3571999999999 -> 369.3334

$ws.Cells.Item(132, 8).Value = 6095.613  # ARM!H132: 6255.567 -> 6095.613
$ws.Cells.Item(132, 9).Value = 1930.1333  # ARM!I132: 1975.3572 -> 1930.1333
$ws.Cells.Item(132, 11).Value = 5790.3999  # ARM!K132: 5926.071599999999 -> 5790.3999
$ws.Cells.Item(132, 13).Value = -3260.3999  # ARM!M132: -3396.071599999999 -> -3260.3999

$ws.Cells.Item(136, 8).Value = 2400  # ARM!H136: 2398.3333 -> 2400
$ws.Cells.Item(136, 9).Value = 0  # ARM!I136: 2397.5 -> 0
$ws.Cells.Item(136, 11).Value = 0  # ARM!K136: 7192.5 -> 0
$ws.Cells.Item(136, 13).ClearContents()  # ARM!M136: -4642.5 -> (removed)

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 5348.6665  # BSM!H20: 5163.8 -> 5348.6665
$ws.Cells.Item(20, 9).Value = 2574.5  # BSM!I20: 2759.6 -> 2574.5
$ws.Cells.Item(20, 11).Value = 2574.5  # BSM!K20: 2759.6 -> 2574.5
$ws.Cells.Item(20, 13).Value = -2327.5  # BSM!M20: -2512.6 -> -2327.5

$ws.Cells.Item(94, 8).Value = 727.25  # BSM!H94: 732.2 -> 727.25
$ws.Cells.Item(94, 9).Value = 686.5789  # BSM!I94: 691.7895 -> 686.5789
$ws.Cells.Item(94, 11).Value = 686.5789  # BSM!K94: 691.7895 -> 686.5789
$ws.Cells.Item(94, 13).Value = -235.5789  # BSM!M94: -240.7895 -> -235.5789

$ws.Cells.Item(99, 8).Value = 2350  # BSM!H99: 1858.25 -> 2350
$ws.Cells.Item(99, 9).Value = 2257.8572  # BSM!I99: 1754.909 -> 2257.8572
$ws.Cells.Item(99, 11).Value = 2257.8572  # BSM!K99: 1754.909 -> 2257.8572
$ws.Cells.Item(99, 13).Value = -759.8571999999999  # BSM!M99: -256.9090000000001 -> -759.8571999999999

$ws.Cells.Item(105, 8).Value = 4084.75  # BSM!H105: 4021.2727 -> 4084.75
$ws.Cells.Item(105, 9).Value = 3283.476  # BSM!I105: 3224.682 -> 3283.476
$ws.Cells.Item(105, 11).Value = 3283.476  # BSM!K105: 3224.682 -> 3283.476
$ws.Cells.Item(105, 13).Value = -1536.476  # BSM!M105: -1477.682 -> -1536.476

$ws.Cells.Item(134, 8).Value = 4982.154  # BSM!H134: 5130.6665 -> 4982.154
$ws.Cells.Item(134, 9).Value = 3696.7778  # BSM!I134: 3758.875 -> 3696.7778
$ws.Cells.Item(134, 11).Value = 11090.3334  # BSM!K134: 11276.625 -> 11090.3334
$ws.Cells.Item(134, 13).Value = -8555.3334  # BSM!M134: -8741.625 -> -8555.3334

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 245  # CRP!H7: 234.125 -> 245
$ws.Cells.Item(7, 9).Value = 268  # CRP!I7: 354.2 -> 268
$ws.Cells.Item(7, 10).Value = 183.66667  # CRP!J7: 34 -> 183.66667
$ws.Cells.Item(7, 11).Value = 268  # CRP!K7: 354.2 -> 268
$ws.Cells.Item(7, 12).Value = 183.66667  # CRP!L7: 34 -> 183.66667
$ws.Cells.Item(7, 13).Value = -155  # CRP!M7: -241.2 -> -155
$ws.Cells.Item(7, 14).Value = -409.66667  # CRP!N7: -260 -> -409.66667

$ws.Cells.Item(16, 8).Value = 369.16666  # CRP!H16: 401.4 -> 369.16666
$ws.Cells.Item(16, 9).Value = 353  # CRP!I16: 389.25 -> 353
$ws.Cells.Item(16, 11).Value = 353  # CRP!K16: 389.25 -> 353
$ws.Cells.Item(16, 13).Value = -66  # CRP!M16: -102.25 -> -66

$ws.Cells.Item(28, 8).Value = 19924.334  # CRP!H28: 21249 -> 19924.334
$ws.Cells.Item(28, 10).Value = 19924.334  # CRP!J28: 21249 -> 19924.334
$ws.Cells.Item(28, 12).Value = 19924.334  # CRP!L28: 21249 -> 19924.334
$ws.Cells.Item(28, 14).Value = -20414.334  # CRP!N28: -21739 -> -20414.334

$ws.Cells.Item(58, 8).Value = 4436.476  # CRP!H58: 4368.7144 -> 4436.476
$ws.Cells.Item(58, 9).Value = 2462.6924  # CRP!I58: 2464.6428 -> 2462.6924
$ws.Cells.Item(58, 10).Value = 7643.875  # CRP!J58: 8176.857 -> 7643.875
$ws.Cells.Item(58, 11).Value = 2462.6924  # CRP!K58: 2464.6428 -> 2462.6924
$ws.Cells.Item(58, 12).Value = 7643.875  # CRP!L58: 8176.857 -> 7643.875
$ws.Cells.Item(58, 13).Value = -2259.6924  # CRP!M58: -2261.6428 -> -2259.6924
$ws.Cells.Item(58, 14).Value = -8049.875  # CRP!N58: -8582.857 -> -8049.875

$ws.Cells.Item(105, 8).Value = 2860.923  # CRP!H105: 2932.8333 -> 2860.923
$ws.Cells.Item(105, 9).Value = 2639.889  # CRP!I105: 2720.125 -> 2639.889
$ws.Cells.Item(105, 11).Value = 2639.889  # CRP!K105: 2720.125 -> 2639.889
$ws.Cells.Item(105, 13).Value = -892.8890000000001  # CRP!M105: -973.125 -> -892.8890000000001

$ws.Cells.Item(113, 8).Value = 369.16666  # CRP!H113: 401.4 -> 369.16666
$ws.Cells.Item(113, 9).Value = 353  # CRP!I113: 389.25 -> 353
$ws.Cells.Item(113, 11).Value = 353  # CRP!K113: 389.25 -> 353
$ws.Cells.Item(113, 13).Value = 1817  # CRP!M113: 1780.75 -> 1817

$ws.Cells.Item(134, 8).Value = 4031.3845  # CRP!H134: 2667.7727 -> 4031.3845
$ws.Cells.Item(134, 9).Value = 3167.2  # CRP!I134: 1997.6316 -> 3167.2
$ws.Cells.Item(134, 11).Value = 9501.599999999999  # CRP!K134: 5992.8948 -> 9501.599999999999
$ws.Cells.Item(134, 13).Value = -6966.599999999999  # CRP!M134: -3457.8948 -> -6966.599999999999

$ws.Cells.Item(136, 8).Value = 4436.476  # CRP!H136: 4368.7144 -> 4436.476
$ws.Cells.Item(136, 9).Value = 2462.6924  # CRP!I136: 2464.6428 -> 2462.6924
$ws.Cells.Item(136, 10).Value = 7643.875  # CRP!J136: 8176.857 -> 7643.875
$ws.Cells.Item(136, 11).Value = 7388.0772  # CRP!K136: 7393.928400000001 -> 7388.0772
$ws.Cells.Item(136, 12).Value = 22931.625  # CRP!L136: 24530.571 -> 22931.625
$ws.Cells.Item(136, 13).Value = -4838.0772  # CRP!M136: -4843.928400000001 -> -4838.0772
$ws.Cells.Item(136, 14).Value = -28031.625  # CRP!N136: -29630.571 -> -28031.625

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(122, 8).Value = 1959.5555  # CUL!H122: 2484.5715 -> 1959.5555
$ws.Cells.Item(122, 10).Value = 2039.5  # CUL!J122: 2998.25 -> 2039.5
$ws.Cells.Item(122, 12).Value = 18355.5  # CUL!L122: 26984.25 -> 18355.5
$ws.Cells.Item(122, 14).Value = -23255.5  # CUL!N122: -31884.25 -> -23255.5

$ws.Cells.Item(132, 8).Value = 2246.5  # CUL!H132: 1998 -> 2246.5
$ws.Cells.Item(132, 10).Value = 2246.5  # CUL!J132: 1998 -> 2246.5
$ws.Cells.Item(132, 12).Value = 20218.5  # CUL!L132: 17982 -> 20218.5
$ws.Cells.Item(132, 14).Value = -25278.5  # CUL!N132: -23042 -> -25278.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 3000  # LTW!H7: 2749.75 -> 3000
$ws.Cells.Item(7, 10).Value = 5000  # LTW!J7: 3499.5 -> 5000
$ws.Cells.Item(7, 12).Value = 5000  # LTW!L7: 3499.5 -> 5000
$ws.Cells.Item(7, 14).Value = -5224  # LTW!N7: -3723.5 -> -5224

$ws.Cells.Item(55, 8).Value = 848.5925999999999  # LTW!H55: 877.4231 -> 848.5925999999999
$ws.Cells.Item(55, 10).Value = 1331.625  # LTW!J55: 1507.7142 -> 1331.625
$ws.Cells.Item(55, 12).Value = 1331.625  # LTW!L55: 1507.7142 -> 1331.625
$ws.Cells.Item(55, 14).Value = -1677.625  # LTW!N55: -1853.7142 -> -1677.625

$ws.Cells.Item(61, 8).Value = 3677.2144  # LTW!H61: 3915.6924 -> 3677.2144
$ws.Cells.Item(61, 9).Value = 3373.625  # LTW!I61: 3627.8635 -> 3373.625
$ws.Cells.Item(61, 11).Value = 3373.625  # LTW!K61: 3627.8635 -> 3373.625
$ws.Cells.Item(61, 13).Value = -3171.625  # LTW!M61: -3425.8635 -> -3171.625

$ws.Cells.Item(93, 8).Value = 936.3333  # LTW!H93: 1046.7 -> 936.3333
$ws.Cells.Item(93, 9).Value = 470.66666  # LTW!I93: 495.2857 -> 470.66666
$ws.Cells.Item(93, 11).Value = 470.66666  # LTW!K93: 495.2857 -> 470.66666
$ws.Cells.Item(93, 13).Value = 777.33334  # LTW!M93: 752.7143 -> 777.33334

$ws.Cells.Item(100, 8).Value = 2119.6  # LTW!H100: 2197.5 -> 2119.6
$ws.Cells.Item(100, 9).Value = 2149.5  # LTW!I100: 2246.875 -> 2149.5
$ws.Cells.Item(100, 11).Value = 2149.5  # LTW!K100: 2246.875 -> 2149.5
$ws.Cells.Item(100, 13).Value = -1608.5  # LTW!M100: -1705.875 -> -1608.5

$ws.Cells.Item(113, 8).Value = 3677.2144  # LTW!H113: 3915.6924 -> 3677.2144
$ws.Cells.Item(113, 9).Value = 3373.625  # LTW!I113: 3627.8635 -> 3373.625
$ws.Cells.Item(113, 11).Value = 3373.625  # LTW!K113: 3627.8635 -> 3373.625
$ws.Cells.Item(113, 13).Value = -1203.625  # LTW!M113: -1457.8635 -> -1203.625

$ws.Cells.Item(126, 8).Value = 3000  # LTW!H126: 2749.75 -> 3000
$ws.Cells.Item(126, 10).Value = 5000  # LTW!J126: 3499.5 -> 5000
$ws.Cells.Item(126, 12).Value = 15000  # LTW!L126: 10498.5 -> 15000
$ws.Cells.Item(126, 14).Value = -19940  # LTW!N126: -15438.5 -> -19940

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(100, 9).Value = 2099.5715  # WVR!I100: 2324.8333 -> 2099.5715
$ws.Cells.Item(100, 10).Value = 1875  # WVR!J100: 1499.3334 -> 1875
$ws.Cells.Item(100, 11).Value = 4199.143  # WVR!K100: 4649.6666 -> 4199.143
$ws.Cells.Item(100, 12).Value = 3750  # WVR!L100: 2998.6668 -> 3750
$ws.Cells.Item(100, 13).Value = -3658.143  # WVR!M100: -4108.6666 -> -3658.143
$ws.Cells.Item(100, 14).Value = -4832  # WVR!N100: -4080.6668 -> -4832

$ws.Cells.Item(136, 8).Value = 87077.586  # WVR!H136: 94539.63 -> 87077.586
$ws.Cells.Item(136, 9).Value = 1863.1428  # WVR!I136: 1863.8572 -> 1863.1428
$ws.Cells.Item(136, 10).Value = 206377.8  # WVR!J136: 256722.25 -> 206377.8
$ws.Cells.Item(136, 11).Value = 5589.428400000001  # WVR!K136: 5591.571599999999 -> 5589.428400000001
$ws.Cells.Item(136, 12).Value = 619133.3999999999  # WVR!L136: 770166.75 -> 619133.3999999999
$ws.Cells.Item(136, 13).Value = -3039.428400000001  # WVR!M136: -3041.571599999999 -> -3039.428400000001
$ws.Cells.Item(136, 14).Value = -624233.3999999999  # WVR!N136: -775266.75 -> -624233.3999999999
